$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'A1'
$ws.Cells.Item(2, 2).Value = 'F0673-GACTAGAACA'
$ws.Cells.Item(2, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGACTAGAACATCGTCGGCAGCGTC'
$ws.Cells.Item(3, 1).Value = 'A2'
$ws.Cells.Item(3, 2).Value = 'F0674-TCTGGATCGT'
$ws.Cells.Item(3, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTCTGGATCGTTCGTCGGCAGCGTC'
$ws.Cells.Item(4, 1).Value = 'A3'
$ws.Cells.Item(4, 2).Value = 'F0675-TCTCACAACA'
$ws.Cells.Item(4, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTCTCACAACATCGTCGGCAGCGTC'
$ws.Cells.Item(5, 1).Value = 'A4'
$ws.Cells.Item(5, 2).Value = 'F0676-GAACTGTGAG'
$ws.Cells.Item(5, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGAACTGTGAGTCGTCGGCAGCGTC'
$ws.Cells.Item(6, 1).Value = 'A5'
$ws.Cells.Item(6, 2).Value = 'F0677-AGTGAGACAC'
$ws.Cells.Item(6, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAGTGAGACACTCGTCGGCAGCGTC'
$ws.Cells.Item(7, 1).Value = 'A6'
$ws.Cells.Item(7, 2).Value = 'F0678-GTAGTGAACA'
$ws.Cells.Item(7, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGTAGTGAACATCGTCGGCAGCGTC'
$ws.Cells.Item(8, 1).Value = 'A7'
$ws.Cells.Item(8, 2).Value = 'F0679-TCTGCACTTC'
$ws.Cells.Item(8, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTCTGCACTTCTCGTCGGCAGCGTC'
$ws.Cells.Item(9, 1).Value = 'A8'
$ws.Cells.Item(9, 2).Value = 'F0680-TGAGGATGCT'
$ws.Cells.Item(9, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTGAGGATGCTTCGTCGGCAGCGTC'
$ws.Cells.Item(10, 1).Value = 'A9'
$ws.Cells.Item(10, 2).Value = 'F0681-AGGATGACTT'
$ws.Cells.Item(10, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAGGATGACTTTCGTCGGCAGCGTC'
$ws.Cells.Item(11, 1).Value = 'A10'
$ws.Cells.Item(11, 2).Value = 'F0682-ACAGCTCTCT'
$ws.Cells.Item(11, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACACAGCTCTCTTCGTCGGCAGCGTC'
$ws.Cells.Item(12, 1).Value = 'A11'
$ws.Cells.Item(12, 2).Value = 'F0683-TGTTGGTAGC'
$ws.Cells.Item(12, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTGTTGGTAGCTCGTCGGCAGCGTC'
$ws.Cells.Item(13, 1).Value = 'A12'
$ws.Cells.Item(13, 2).Value = 'F0684-TACGACTAGT'
$ws.Cells.Item(13, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTACGACTAGTTCGTCGGCAGCGTC'
$ws.Cells.Item(14, 1).Value = 'B1'
$ws.Cells.Item(14, 2).Value = 'F0685-ACATGGAGTA'
$ws.Cells.Item(14, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACACATGGAGTATCGTCGGCAGCGTC'
$ws.Cells.Item(15, 1).Value = 'B2'
$ws.Cells.Item(15, 2).Value = 'F0686-TAGACAACGA'
$ws.Cells.Item(15, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTAGACAACGATCGTCGGCAGCGTC'
$ws.Cells.Item(16, 1).Value = 'B3'
$ws.Cells.Item(16, 2).Value = 'F0687-ATGATCAGCA'
$ws.Cells.Item(16, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACATGATCAGCATCGTCGGCAGCGTC'
$ws.Cells.Item(17, 1).Value = 'B4'
$ws.Cells.Item(17, 2).Value = 'F0688-ATGAACGAGA'
$ws.Cells.Item(17, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACATGAACGAGATCGTCGGCAGCGTC'
$ws.Cells.Item(18, 1).Value = 'B5'
$ws.Cells.Item(18, 2).Value = 'F0689-CAAGTGACGT'
$ws.Cells.Item(18, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCAAGTGACGTTCGTCGGCAGCGTC'
$ws.Cells.Item(19, 1).Value = 'B6'
$ws.Cells.Item(19, 2).Value = 'F0690-AGTACAAGTC'
$ws.Cells.Item(19, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAGTACAAGTCTCGTCGGCAGCGTC'
$ws.Cells.Item(20, 1).Value = 'B7'
$ws.Cells.Item(20, 2).Value = 'F0691-ACACAGTGTA'
$ws.Cells.Item(20, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACACACAGTGTATCGTCGGCAGCGTC'
$ws.Cells.Item(21, 1).Value = 'B8'
$ws.Cells.Item(21, 2).Value = 'F0692-ACGTGACGTA'
$ws.Cells.Item(21, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACACGTGACGTATCGTCGGCAGCGTC'
$ws.Cells.Item(22, 1).Value = 'B9'
$ws.Cells.Item(22, 2).Value = 'F0693-AAGTACGAGG'
$ws.Cells.Item(22, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAAGTACGAGGTCGTCGGCAGCGTC'
$ws.Cells.Item(23, 1).Value = 'B10'
$ws.Cells.Item(23, 2).Value = 'F0694-AGGTACCTCA'
$ws.Cells.Item(23, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAGGTACCTCATCGTCGGCAGCGTC'
$ws.Cells.Item(24, 1).Value = 'B11'
$ws.Cells.Item(24, 2).Value = 'F0695-ACAACGAAGG'
$ws.Cells.Item(24, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACACAACGAAGGTCGTCGGCAGCGTC'
$ws.Cells.Item(25, 1).Value = 'B12'
$ws.Cells.Item(25, 2).Value = 'F0696-AGCATGCTAG'
$ws.Cells.Item(25, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAGCATGCTAGTCGTCGGCAGCGTC'
$ws.Cells.Item(26, 1).Value = 'C1'
$ws.Cells.Item(26, 2).Value = 'F0697-CACAGAGTGT'
$ws.Cells.Item(26, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCACAGAGTGTTCGTCGGCAGCGTC'
$ws.Cells.Item(27, 1).Value = 'C2'
$ws.Cells.Item(27, 2).Value = 'F0698-CAAGACTAGG'
$ws.Cells.Item(27, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCAAGACTAGGTCGTCGGCAGCGTC'
$ws.Cells.Item(28, 1).Value = 'C3'
$ws.Cells.Item(28, 2).Value = 'F0699-AGTTCGTGAG'
$ws.Cells.Item(28, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAGTTCGTGAGTCGTCGGCAGCGTC'
$ws.Cells.Item(29, 1).Value = 'C4'
$ws.Cells.Item(29, 2).Value = 'F0700-TCACTTCAAG'
$ws.Cells.Item(29, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTCACTTCAAGTCGTCGGCAGCGTC'
$ws.Cells.Item(30, 1).Value = 'C5'
$ws.Cells.Item(30, 2).Value = 'F0701-TCTGCTGTAG'
$ws.Cells.Item(30, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTCTGCTGTAGTCGTCGGCAGCGTC'
$ws.Cells.Item(31, 1).Value = 'C6'
$ws.Cells.Item(31, 2).Value = 'F0702-AACAAGACCA'
$ws.Cells.Item(31, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAACAAGACCATCGTCGGCAGCGTC'
$ws.Cells.Item(32, 1).Value = 'C7'
$ws.Cells.Item(32, 2).Value = 'F0703-ACTGACTCAT'
$ws.Cells.Item(32, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACACTGACTCATTCGTCGGCAGCGTC'
$ws.Cells.Item(33, 1).Value = 'C8'
$ws.Cells.Item(33, 2).Value = 'F0704-AACATCCTGA'
$ws.Cells.Item(33, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAACATCCTGATCGTCGGCAGCGTC'
$ws.Cells.Item(34, 1).Value = 'C9'
$ws.Cells.Item(34, 2).Value = 'F0705-GATCACGAAC'
$ws.Cells.Item(34, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGATCACGAACTCGTCGGCAGCGTC'
$ws.Cells.Item(35, 1).Value = 'C10'
$ws.Cells.Item(35, 2).Value = 'F0706-CTGTTCTAGC'
$ws.Cells.Item(35, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCTGTTCTAGCTCGTCGGCAGCGTC'
$ws.Cells.Item(36, 1).Value = 'C11'
$ws.Cells.Item(36, 2).Value = 'F0707-ACTAGAGACC'
$ws.Cells.Item(36, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACACTAGAGACCTCGTCGGCAGCGTC'
$ws.Cells.Item(37, 1).Value = 'C12'
$ws.Cells.Item(37, 2).Value = 'F0708-GCTCTGAGAA'
$ws.Cells.Item(37, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGCTCTGAGAATCGTCGGCAGCGTC'
$ws.Cells.Item(38, 1).Value = 'D1'
$ws.Cells.Item(38, 2).Value = 'F0709-TAGACACCTG'
$ws.Cells.Item(38, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTAGACACCTGTCGTCGGCAGCGTC'
$ws.Cells.Item(39, 1).Value = 'D2'
$ws.Cells.Item(39, 2).Value = 'F0710-AACGTAGTCT'
$ws.Cells.Item(39, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAACGTAGTCTTCGTCGGCAGCGTC'
$ws.Cells.Item(40, 1).Value = 'D3'
$ws.Cells.Item(40, 2).Value = 'F0711-AACCACGAGA'
$ws.Cells.Item(40, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAACCACGAGATCGTCGGCAGCGTC'
$ws.Cells.Item(41, 1).Value = 'D4'
$ws.Cells.Item(41, 2).Value = 'F0712-AGGATGTTGT'
$ws.Cells.Item(41, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAGGATGTTGTTCGTCGGCAGCGTC'
$ws.Cells.Item(42, 1).Value = 'D5'
$ws.Cells.Item(42, 2).Value = 'F0713-ACGAGTCTGA'
$ws.Cells.Item(42, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACACGAGTCTGATCGTCGGCAGCGTC'
$ws.Cells.Item(43, 1).Value = 'D6'
$ws.Cells.Item(43, 2).Value = 'F0714-TGTAGAGTGC'
$ws.Cells.Item(43, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTGTAGAGTGCTCGTCGGCAGCGTC'
$ws.Cells.Item(44, 1).Value = 'D7'
$ws.Cells.Item(44, 2).Value = 'F0715-GTTGCATCTC'
$ws.Cells.Item(44, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGTTGCATCTCTCGTCGGCAGCGTC'
$ws.Cells.Item(45, 1).Value = 'D8'
$ws.Cells.Item(45, 2).Value = 'F0716-TGCTAGCACA'
$ws.Cells.Item(45, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTGCTAGCACATCGTCGGCAGCGTC'
$ws.Cells.Item(46, 1).Value = 'D9'
$ws.Cells.Item(46, 2).Value = 'F0717-GACTGCTCAT'
$ws.Cells.Item(46, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGACTGCTCATTCGTCGGCAGCGTC'
$ws.Cells.Item(47, 1).Value = 'D10'
$ws.Cells.Item(47, 2).Value = 'F0718-CAGGATGTCA'
$ws.Cells.Item(47, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCAGGATGTCATCGTCGGCAGCGTC'
$ws.Cells.Item(48, 1).Value = 'D11'
$ws.Cells.Item(48, 2).Value = 'F0719-TCACCATCCT'
$ws.Cells.Item(48, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTCACCATCCTTCGTCGGCAGCGTC'
$ws.Cells.Item(49, 1).Value = 'D12'
$ws.Cells.Item(49, 2).Value = 'F0720-CGTGTAGCTT'
$ws.Cells.Item(49, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCGTGTAGCTTTCGTCGGCAGCGTC'
$ws.Cells.Item(50, 1).Value = 'E1'
$ws.Cells.Item(50, 2).Value = 'F0721-CGTGATGATG'
$ws.Cells.Item(50, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCGTGATGATGTCGTCGGCAGCGTC'
$ws.Cells.Item(51, 1).Value = 'E2'
$ws.Cells.Item(51, 2).Value = 'F0722-TAGCAAGACT'
$ws.Cells.Item(51, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTAGCAAGACTTCGTCGGCAGCGTC'
$ws.Cells.Item(52, 1).Value = 'E3'
$ws.Cells.Item(52, 2).Value = 'F0723-GAACCACAGT'
$ws.Cells.Item(52, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGAACCACAGTTCGTCGGCAGCGTC'
$ws.Cells.Item(53, 1).Value = 'E4'
$ws.Cells.Item(53, 2).Value = 'F0724-TCAACTGGAT'
$ws.Cells.Item(53, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTCAACTGGATTCGTCGGCAGCGTC'
$ws.Cells.Item(54, 1).Value = 'E5'
$ws.Cells.Item(54, 2).Value = 'F0725-GTCAGTACGA'
$ws.Cells.Item(54, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGTCAGTACGATCGTCGGCAGCGTC'
$ws.Cells.Item(55, 1).Value = 'E6'
$ws.Cells.Item(55, 2).Value = 'F0726-TGATGACATC'
$ws.Cells.Item(55, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTGATGACATCTCGTCGGCAGCGTC'
$ws.Cells.Item(56, 1).Value = 'E7'
$ws.Cells.Item(56, 2).Value = 'F0727-GACTCCTTGT'
$ws.Cells.Item(56, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGACTCCTTGTTCGTCGGCAGCGTC'
$ws.Cells.Item(57, 1).Value = 'E8'
$ws.Cells.Item(57, 2).Value = 'F0728-GAACGTACGA'
$ws.Cells.Item(57, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGAACGTACGATCGTCGGCAGCGTC'
$ws.Cells.Item(58, 1).Value = 'E9'
$ws.Cells.Item(58, 2).Value = 'F0729-GCTCAAGCTT'
$ws.Cells.Item(58, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGCTCAAGCTTTCGTCGGCAGCGTC'
$ws.Cells.Item(59, 1).Value = 'E10'
$ws.Cells.Item(59, 2).Value = 'F0730-TTCCTGCTTC'
$ws.Cells.Item(59, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTTCCTGCTTCTCGTCGGCAGCGTC'
$ws.Cells.Item(60, 1).Value = 'E11'
$ws.Cells.Item(60, 2).Value = 'F0731-TCTTGGATGC'
$ws.Cells.Item(60, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTCTTGGATGCTCGTCGGCAGCGTC'
$ws.Cells.Item(61, 1).Value = 'E12'
$ws.Cells.Item(61, 2).Value = 'F0732-CTCTCACTGT'
$ws.Cells.Item(61, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCTCTCACTGTTCGTCGGCAGCGTC'
$ws.Cells.Item(62, 1).Value = 'F1'
$ws.Cells.Item(62, 2).Value = 'F0733-TGATGCTCAT'
$ws.Cells.Item(62, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTGATGCTCATTCGTCGGCAGCGTC'
$ws.Cells.Item(63, 1).Value = 'F2'
$ws.Cells.Item(63, 2).Value = 'F0734-CTCCTCTCAT'
$ws.Cells.Item(63, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCTCCTCTCATTCGTCGGCAGCGTC'
$ws.Cells.Item(64, 1).Value = 'F3'
$ws.Cells.Item(64, 2).Value = 'F0735-TTCCAACGAG'
$ws.Cells.Item(64, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTTCCAACGAGTCGTCGGCAGCGTC'
$ws.Cells.Item(65, 1).Value = 'F4'
$ws.Cells.Item(65, 2).Value = 'F0736-AACATCAGCT'
$ws.Cells.Item(65, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAACATCAGCTTCGTCGGCAGCGTC'
$ws.Cells.Item(66, 1).Value = 'F5'
$ws.Cells.Item(66, 2).Value = 'F0737-AGATCCAACT'
$ws.Cells.Item(66, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAGATCCAACTTCGTCGGCAGCGTC'
$ws.Cells.Item(67, 1).Value = 'F6'
$ws.Cells.Item(67, 2).Value = 'F0738-GAAGAACAGG'
$ws.Cells.Item(67, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGAAGAACAGGTCGTCGGCAGCGTC'
$ws.Cells.Item(68, 1).Value = 'F7'
$ws.Cells.Item(68, 2).Value = 'F0739-TGTGAGCTCT'
$ws.Cells.Item(68, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTGTGAGCTCTTCGTCGGCAGCGTC'
$ws.Cells.Item(69, 1).Value = 'F8'
$ws.Cells.Item(69, 2).Value = 'F0740-GTCGATCTGT'
$ws.Cells.Item(69, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGTCGATCTGTTCGTCGGCAGCGTC'
$ws.Cells.Item(70, 1).Value = 'F9'
$ws.Cells.Item(70, 2).Value = 'F0741-GTACACAGAC'
$ws.Cells.Item(70, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGTACACAGACTCGTCGGCAGCGTC'
$ws.Cells.Item(71, 1).Value = 'F10'
$ws.Cells.Item(71, 2).Value = 'F0742-GTAGCTCAGA'
$ws.Cells.Item(71, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGTAGCTCAGATCGTCGGCAGCGTC'
$ws.Cells.Item(72, 1).Value = 'F11'
$ws.Cells.Item(72, 2).Value = 'F0743-GTCGTCAGTA'
$ws.Cells.Item(72, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGTCGTCAGTATCGTCGGCAGCGTC'
$ws.Cells.Item(73, 1).Value = 'F12'
$ws.Cells.Item(73, 2).Value = 'F0744-GAGTTCTCTT'
$ws.Cells.Item(73, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGAGTTCTCTTTCGTCGGCAGCGTC'
$ws.Cells.Item(74, 1).Value = 'G1'
$ws.Cells.Item(74, 2).Value = 'F0745-AACTCAGAGC'
$ws.Cells.Item(74, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAACTCAGAGCTCGTCGGCAGCGTC'
$ws.Cells.Item(75, 1).Value = 'G2'
$ws.Cells.Item(75, 2).Value = 'F0746-CATGAAGTGG'
$ws.Cells.Item(75, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCATGAAGTGGTCGTCGGCAGCGTC'
$ws.Cells.Item(76, 1).Value = 'G3'
$ws.Cells.Item(76, 2).Value = 'F0747-TGCTCTCCAT'
$ws.Cells.Item(76, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTGCTCTCCATTCGTCGGCAGCGTC'
$ws.Cells.Item(77, 1).Value = 'G4'
$ws.Cells.Item(77, 2).Value = 'F0748-ACCTCTGATG'
$ws.Cells.Item(77, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACACCTCTGATGTCGTCGGCAGCGTC'
$ws.Cells.Item(78, 1).Value = 'G5'
$ws.Cells.Item(78, 2).Value = 'F0749-GGAACCAAGT'
$ws.Cells.Item(78, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGGAACCAAGTTCGTCGGCAGCGTC'
$ws.Cells.Item(79, 1).Value = 'G6'
$ws.Cells.Item(79, 2).Value = 'F0750-TACACTGGTT'
$ws.Cells.Item(79, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTACACTGGTTTCGTCGGCAGCGTC'
$ws.Cells.Item(80, 1).Value = 'G7'
$ws.Cells.Item(80, 2).Value = 'F0751-ATCTGGATCC'
$ws.Cells.Item(80, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACATCTGGATCCTCGTCGGCAGCGTC'
$ws.Cells.Item(81, 1).Value = 'G8'
$ws.Cells.Item(81, 2).Value = 'F0752-ATGACTGTGC'
$ws.Cells.Item(81, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACATGACTGTGCTCGTCGGCAGCGTC'
$ws.Cells.Item(82, 1).Value = 'G9'
$ws.Cells.Item(82, 2).Value = 'F0753-GCTTGAGTAG'
$ws.Cells.Item(82, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGCTTGAGTAGTCGTCGGCAGCGTC'
$ws.Cells.Item(83, 1).Value = 'G10'
$ws.Cells.Item(83, 2).Value = 'F0754-CTGACTAGGA'
$ws.Cells.Item(83, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCTGACTAGGATCGTCGGCAGCGTC'
$ws.Cells.Item(84, 1).Value = 'G11'
$ws.Cells.Item(84, 2).Value = 'F0755-CTTGACGTCA'
$ws.Cells.Item(84, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCTTGACGTCATCGTCGGCAGCGTC'
$ws.Cells.Item(85, 1).Value = 'G12'
$ws.Cells.Item(85, 2).Value = 'F0756-GTGTCTACTC'
$ws.Cells.Item(85, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGTGTCTACTCTCGTCGGCAGCGTC'
$ws.Cells.Item(86, 1).Value = 'H1'
$ws.Cells.Item(86, 2).Value = 'F0757-TTCTACGTGG'
$ws.Cells.Item(86, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACTTCTACGTGGTCGTCGGCAGCGTC'
$ws.Cells.Item(87, 1).Value = 'H2'
$ws.Cells.Item(87, 2).Value = 'F0758-ACTGTGTCAC'
$ws.Cells.Item(87, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACACTGTGTCACTCGTCGGCAGCGTC'
$ws.Cells.Item(88, 1).Value = 'H3'
$ws.Cells.Item(88, 2).Value = 'F0759-GTGAACATGT'
$ws.Cells.Item(88, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGTGAACATGTTCGTCGGCAGCGTC'
$ws.Cells.Item(89, 1).Value = 'H4'
$ws.Cells.Item(89, 2).Value = 'F0760-AGCTTCCTTC'
$ws.Cells.Item(89, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAGCTTCCTTCTCGTCGGCAGCGTC'
$ws.Cells.Item(90, 1).Value = 'H5'
$ws.Cells.Item(90, 2).Value = 'F0761-ACTGCATGGA'
$ws.Cells.Item(90, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACACTGCATGGATCGTCGGCAGCGTC'
$ws.Cells.Item(91, 1).Value = 'H6'
$ws.Cells.Item(91, 2).Value = 'F0762-AAGTGTTCTC'
$ws.Cells.Item(91, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACAAGTGTTCTCTCGTCGGCAGCGTC'
$ws.Cells.Item(92, 1).Value = 'H7'
$ws.Cells.Item(92, 2).Value = 'F0763-GCTGACTTCT'
$ws.Cells.Item(92, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGCTGACTTCTTCGTCGGCAGCGTC'
$ws.Cells.Item(93, 1).Value = 'H8'
$ws.Cells.Item(93, 2).Value = 'F0764-GGTGAGAAGA'
$ws.Cells.Item(93, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGGTGAGAAGATCGTCGGCAGCGTC'
$ws.Cells.Item(94, 1).Value = 'H9'
$ws.Cells.Item(94, 2).Value = 'F0765-GGACTGAACT'
$ws.Cells.Item(94, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGGACTGAACTTCGTCGGCAGCGTC'
$ws.Cells.Item(95, 1).Value = 'H10'
$ws.Cells.Item(95, 2).Value = 'F0766-CAACGTGTAG'
$ws.Cells.Item(95, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACCAACGTGTAGTCGTCGGCAGCGTC'
$ws.Cells.Item(96, 1).Value = 'H11'
$ws.Cells.Item(96, 2).Value = 'F0767-ACCAACTAGA'
$ws.Cells.Item(96, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACACCAACTAGATCGTCGGCAGCGTC'
$ws.Cells.Item(97, 1).Value = 'H12'
$ws.Cells.Item(97, 2).Value = 'F0768-GTGAGTTCTG'
$ws.Cells.Item(97, 3).Value = 'AATGATACGGCGACCACCGAGATCTACACGTGAGTTCTGTCGTCGGCAGCGTC'
